$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"

# Column D: date, copy formatting/style from the cell above (D23) and set the date serial value
$ws.Cells.Item($row, 4).Value = 44448
$ws.Cells.Item(23, 4).Copy()
$ws.Cells.Item($row, 4).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item($row, 4).Value = 44448

$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100114002
$ws.Cells.Item($row, 7).Value = "Camote"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 45
$ws.Cells.Item($row, 11).Value = 20000
$ws.Cells.Item($row, 12).Value = 20000
$ws.Cells.Item($row, 13).Value = 20000
$ws.Cells.Item($row, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item($row, 15).Value = "Perú"
$ws.Cells.Item($row, 16).Value = 1000
$ws.Cells.Item($row, 17).Value = 20
$ws.Cells.Item($row, 18).Value = "Hortaliza"
